# OAB : Laporan Riwayat Pengobatan LUTS sebelumnya
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new blank row at row 6 (pushes the old header-detail row 6 -> 7,
#    and expands the A5:A6 merge to A5:A7 automatically).
# ---------------------------------------------------------------------------
$ws.Rows("6:6").Insert()

# Row 5 becomes a taller header row.
$ws.Rows(5).RowHeight = 75

# ---------------------------------------------------------------------------
# 2) New top-level header cell BR5 ("Riwayat Pengobatan LUTS sebelumnya"),
#    merged across BR5:CT5, styled like the other top-level group headers.
# ---------------------------------------------------------------------------
$top = $ws.Range("BR5:CT5")
$top.Merge()
$top.Value = "Riwayat Pengobatan LUTS sebelumnya"
$top.Font.Bold = $true
$top.HorizontalAlignment = -4131   # xlLeft
$top.VerticalAlignment = -4160     # xlTop
$top.WrapText = $true
$top.Interior.ThemeColor = 4
$top.Interior.TintAndShade = 0.79998168889431442

Write-Host "stage1 done"
